$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Cell C6 ("Obiettivi prefissati:" ...): replace the closing sentence
#    about "Sistemi e Reti" with the updated wording, then restore the
#    two bold runs (a line break and a single space) that the text
#    replacement flattens.
# ---------------------------------------------------------------------
$rC6 = $ws.Range("C6")
$fullC6 = $rC6.Value2
$marker = "Considerando invece"
$idx = $fullC6.IndexOf($marker)
$startPos = $idx + 1
$length = $fullC6.Length - $idx

$newTail = "Considerando invece il progetto relativo a Sistemi e Reti il team ha fissato come obiettivo quello di realizzare un infrastruttura informatica tra tre sedi che garantisca uno scambio di dati funzionale,più sicuro e all'avanguardia possibile."

$rC6.Characters($startPos, $length).Text = $newTail

# Restore bold formatting on the line break after the first paragraph
# and on the single space between "dinamica" and "impaginata".
$rC6.Characters(446, 1).Font.Bold = $true
$rC6.Characters(590, 1).Font.Bold = $true

# ---------------------------------------------------------------------
# 2) Cell C9 ("Realizzazione:" ...): append the sentence about the
#    Database/Packet Tracer deliverables, then restore the bold run on
#    the "Realizzazione:" heading that the edit flattens.
# ---------------------------------------------------------------------
$rC9 = $ws.Range("C9")
$fullC9 = $rC9.Value2
$appendText = " Verranno realizzati i due Database per Tp ed Informatica ed il file Packet Tracer contentente la struttura della rete realizzata."
$insertPos = $fullC9.Length + 1

$rC9.Characters($insertPos, 0).Text = $appendText
$rC9.Characters(1, 14).Font.Bold = $true

# ---------------------------------------------------------------------
# 3) Row 9 height: 75 -> 108
# ---------------------------------------------------------------------
$ws.Rows.Item(9).RowHeight = 108

# ---------------------------------------------------------------------
# 4) Selection moves from G9 to C9
# ---------------------------------------------------------------------
$rC9.Select() | Out-Null
